$wb = $excel.ActiveWorkbook

# --- MapSetting sheet: add the new NPC row (row 5) ---------------------
$map = $wb.Worksheets.Item("MapSetting")

# Order matters for the shared-strings table: write the new unique strings
# in the same order they first appear left-to-right on the new row so the
# generated <sst> indices line up (npcsainisi=79, 塞尼斯=80).
$map.Range("D5").Value = "npcsainisi"
$map.Range("E5").Value = "npcsainisi"
$map.Range("B5").Value = "塞尼斯"

$map.Range("A5").Value = 42100002
$map.Range("C5").Value = 0
# Leading apostrophe forces text (otherwise "true" is auto-typed as Boolean)
$map.Range("F5").Value = "'true"
$map.Range("G5").Value = 43000101
$map.Range("H5").Value = "oneline"
$map.Range("N5").Value = 22036101

# Match the style used on the row above (borders on H, number format on F)
$map.Range("H5").Style = $map.Range("H4").Style
$map.Range("F5").Style = $map.Range("F4").Style

# Grow the worksheet table ("表1_3") so it covers the new row too
$t = $map.ListObjects.Item(1)
$t.Resize($map.Range("A1:S5"))

# Update selection to match the authored state
$map.Range("K4").Select()

# --- SceneQuest sheet: just move the selection -------------------------
$scene = $wb.Worksheets.Item("SceneQuest")
$scene.Range("H6").Select()

# Re-activate MapSetting as the active sheet/tab
$map.Activate()
